$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("L4").Value = 82
